$wb = $excel.ActiveWorkbook

# Sheet "OFF" - row 3 (R row)
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 233
$wsOff.Range("C3").Value = 170
$wsOff.Range("D3").Value = 42
$wsOff.Range("E3").Value = 22
$wsOff.Range("F3").Value = 2
$wsOff.Range("G3").Value = 5

# Sheet "DEF" - row 3 (R row)
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 161
$wsDef.Range("C3").Value = 117
$wsDef.Range("D3").Value = 42
$wsDef.Range("E3").Value = 22
